$d = $word.ActiveDocument

# Helper: XML package wrapper for InsertXML fragments targeting word/document.xml
function New-DocXmlPackage([string]$bodyXml) {
    return '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + $bodyXml + '</w:document></pkg:xmlData></pkg:part></pkg:package>'
}

# ---------------------------------------------------------------------------
# 1. Paragraph 1: wrap the existing runs with a "_GoBack" bookmark
#    ("FRDR Back Plan for Geodisy" + " in case of a full crash")
# ---------------------------------------------------------------------------
$p1 = $d.Paragraphs(1).Range
$r1 = $d.Range($p1.Start, $p1.End)
$d.Bookmarks.Add("_GoBack", $r1) | Out-Null

# ---------------------------------------------------------------------------
# 2. Paragraph 5 ("Metadata (...) in /var/www/[webserver]/html/geodisy/")
#    - merge "Metadata (" + "Geoblacklight" + " json and ISO XML)" into one run
#    - merge "]/html/" + "geodisy" + "/" into one run
#    (drop the proofErr wrappers, text unchanged)
# ---------------------------------------------------------------------------
$body5 = '<w:body><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr>' + `
    '<w:r><w:t>Metadata (Geoblacklight json and ISO XML)</w:t></w:r>' + `
    '<w:r><w:t xml:space="preserve"> in /var/www/[</w:t></w:r>' + `
    '<w:r><w:t>webserver</w:t></w:r>' + `
    '<w:r><w:t>]/html/geodisy/</w:t></w:r>' + `
    '</w:p></w:body>'
$d.Paragraphs(5).Range.InsertXML((New-DocXmlPackage $body5)) | Out-Null

# ---------------------------------------------------------------------------
# 3. Paragraph 6 ("Solr index for GeoBlacklight")
#    - merge "Solr" + " index for GeoBlacklight" into one run (drop proofErr)
# ---------------------------------------------------------------------------
$body6 = '<w:body><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr>' + `
    '<w:r><w:t>Solr index for GeoBlacklight</w:t></w:r>' + `
    '</w:p></w:body>'
$d.Paragraphs(6).Range.InsertXML((New-DocXmlPackage $body6)) | Out-Null

# ---------------------------------------------------------------------------
# 4. Paragraph 7 ("Password file (which is not stored in GitHub for security reasons)")
#    - merge "Password " + "file " into "Password file" (note: trailing space
#      from "file " is dropped, the lone single-space run right after supplies
#      the only space between "file" and "(")
#    - merge "(" + "which is not stored in " into "(which is not stored in "
#    - fix "Github" -> "GitHub" as its own run
#    - keep trailing " for security reasons)"
#    (drop the gramStart/gramEnd/spellStart/spellEnd proofErr wrappers)
# ---------------------------------------------------------------------------
$body7 = '<w:body><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr>' + `
    '<w:r><w:t>Password file</w:t></w:r>' + `
    '<w:r><w:t xml:space="preserve"> </w:t></w:r>' + `
    '<w:r><w:t xml:space="preserve">(which is not stored in </w:t></w:r>' + `
    '<w:r><w:t>GitHub</w:t></w:r>' + `
    '<w:r><w:t xml:space="preserve"> for security reasons)</w:t></w:r>' + `
    '</w:p></w:body>'
$d.Paragraphs(7).Range.InsertXML((New-DocXmlPackage $body7)) | Out-Null

# ---------------------------------------------------------------------------
# 5. Paragraph 9 ("All files in /home/centos/Geodisy/savedFiles/")
#    - merge "/home/centos/Geodisy/" + "savedFiles" into one run (drop proofErr)
# ---------------------------------------------------------------------------
$body9 = '<w:body><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr>' + `
    '<w:r><w:t>All</w:t></w:r>' + `
    '<w:r><w:t xml:space="preserve"> files</w:t></w:r>' + `
    '<w:r><w:t xml:space="preserve"> in</w:t></w:r>' + `
    '<w:r><w:t xml:space="preserve"> </w:t></w:r>' + `
    '<w:r><w:t>/home/centos/Geodisy/savedFiles</w:t></w:r>' + `
    '<w:r><w:t>/</w:t></w:r>' + `
    '</w:p></w:body>'
$d.Paragraphs(9).Range.InsertXML((New-DocXmlPackage $body9)) | Out-Null

# ---------------------------------------------------------------------------
# 6. Paragraph 10 ("All files in /opt/geoblacklight/geodisy/")
#    - remove the "_GoBack" bookmark that used to sit here (it moved to para 1)
# ---------------------------------------------------------------------------
$body10 = '<w:body><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr>' + `
    '<w:r><w:t xml:space="preserve">All files in </w:t></w:r>' + `
    '<w:r><w:t>/opt/geoblacklight/geodisy/</w:t></w:r>' + `
    '</w:p></w:body>'
$d.Paragraphs(10).Range.InsertXML((New-DocXmlPackage $body10)) | Out-Null

# ---------------------------------------------------------------------------
# 7. Paragraph 11 ("Once Geoserver has been implemented")
#    - merge "Once " + "Geoserver" + " has been implemented" into one run (drop proofErr)
# ---------------------------------------------------------------------------
$body11 = '<w:body><w:p><w:r><w:t>Once Geoserver has been implemented</w:t></w:r></w:p></w:body>'
$d.Paragraphs(11).Range.InsertXML((New-DocXmlPackage $body11)) | Out-Null

# ---------------------------------------------------------------------------
# 8. Paragraph 12 ("PostGis Database" -> "PostGIS Database")
#    - fix "PostGis" -> "PostGIS" split across two runs "PostG" + "IS"
#    - drop the spellStart/spellEnd proofErr wrapper
# ---------------------------------------------------------------------------
$body12 = '<w:body><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr>' + `
    '<w:r><w:t>PostG</w:t></w:r>' + `
    '<w:r><w:t>IS</w:t></w:r>' + `
    '<w:r><w:t xml:space="preserve"> Database</w:t></w:r>' + `
    '</w:p></w:body>'
$d.Paragraphs(12).Range.InsertXML((New-DocXmlPackage $body12)) | Out-Null
